$d = $word.ActiveDocument

# --- Step 1: delete the "License Information" Heading2 paragraph (was paragraph 4) ---
$pLicenseInfo = $d.Paragraphs.Item(4)
$pLicenseInfo.Range.Delete()

# --- Step 2: rebuild the license-notice paragraph (now paragraph 4) ---
$pLic = $d.Paragraphs.Item(4)
$innerStart = $pLic.Range.Start
$innerEnd = $pLic.Range.End - 1
$clearRng = $d.Range($innerStart, $innerEnd)
$clearRng.Text = ""

$pos = $pLic.Range.Start

$seg1 = "Biblica Study Notes (Key Terms)"
$r1 = $d.Range($pos, $pos)
$r1.InsertAfter($seg1)
$r1.Font.Bold = 1
$pos = $pos + $seg1.Length

$seg2 = " © 2023 Biblica Inc. Released under CC BY-SA 4.0 license. "
$r2 = $d.Range($pos, $pos)
$r2.InsertAfter($seg2)
$r2.Font.Bold = 0
$pos = $pos + $seg2.Length

$seg3 = "Biblica Study Notes"
$r3 = $d.Range($pos, $pos)
$r3.InsertAfter($seg3)
$r3.Font.Bold = 0
$pos = $pos + $seg3.Length

$seg4 = " has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文)from Biblica Study Notes © 2023 Biblica Inc. Released under CC BY-SA 4.0 license by Mission Mutual."
$r4 = $d.Range($pos, $pos)
$r4.InsertAfter($seg4)
$r4.Font.Bold = 0
$pos = $pos + $seg4.Length

# --- Step 3: delete the "This PDF version is provided under the same license." paragraph (now paragraph 5) ---
$pPdf = $d.Paragraphs.Item(5)
$pPdf.Range.Delete()

# --- Step 4: delete the italic "व्यवस्था के शिक्षक" body paragraph that follows the "वय" Heading2 ---
$pVay = $d.Paragraphs.Item(10)
$pItalic = $d.Paragraphs.Item(11)
$pItalic.Range.Delete()

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
